# Auto-generated: apply value changes described by the xml_diff.
# Grouped per worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW); WVR is untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 299
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H17").Value = 2738.1875
$ws.Range("J17").Value = 2349.8667
$ws.Range("L17").Value = 7049.6001
$ws.Range("N17").Value = -7385.6001
$ws.Range("H28").Value = 233.72728
$ws.Range("I28").Value = 157.1
$ws.Range("K28").Value = 157.1
$ws.Range("M28").Value = 327.9
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H70").Value = 18474.875
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 20828.428
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 62485.284
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -63025.284
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H73").Value = 18474.875
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 20828.428
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 62485.284
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -64357.284
$ws.Range("H82").Value = 3500
$ws.Range("I82").Value = 3500
$ws.Range("K82").Value = 10500
$ws.Range("M82").Value = -10094
$ws.Range("H85").Value = 3500
$ws.Range("I85").Value = 3500
$ws.Range("K85").Value = 10500
$ws.Range("M85").Value = -9096
$ws.Range("H86").Value = 1739.4
$ws.Range("I86").Value = 1749.5
$ws.Range("J86").Value = 1732.6666
$ws.Range("K86").Value = 1749.5
$ws.Range("L86").Value = 1732.6666
$ws.Range("M86").Value = -626.5
$ws.Range("N86").Value = -3978.6666
$ws.Range("H88").Value = 35716036
$ws.Range("J88").Value = 1966.3334
$ws.Range("L88").Value = 1966.3334
$ws.Range("N88").Value = -2778.3334
$ws.Range("H89").Value = 1739.4
$ws.Range("I89").Value = 1749.5
$ws.Range("J89").Value = 1732.6666
$ws.Range("K89").Value = 8747.5
$ws.Range("L89").Value = 8663.333
$ws.Range("M89").Value = -3131.5
$ws.Range("N89").Value = -19895.333
$ws.Range("H91").Value = 35716036
$ws.Range("J91").Value = 1966.3334
$ws.Range("L91").Value = 1966.3334
$ws.Range("N91").Value = -4774.3334
$ws.Range("H98").Value = 864.84375
$ws.Range("I98").Value = 570.4583
$ws.Range("J98").Value = 1748
$ws.Range("K98").Value = 570.4583
$ws.Range("L98").Value = 1748
$ws.Range("M98").Value = 927.5417
$ws.Range("N98").Value = -4744
$ws.Range("H112").Value = 1424.1786
$ws.Range("J112").Value = 1424.1786
$ws.Range("L112").Value = 4272.5358
$ws.Range("N112").Value = -6488.5358
$ws.Range("H122").Value = 864.84375
$ws.Range("I122").Value = 570.4583
$ws.Range("J122").Value = 1748
$ws.Range("K122").Value = 1711.3749
$ws.Range("L122").Value = 5244
$ws.Range("M122").Value = 738.6251
$ws.Range("N122").Value = -10144
$ws.Range("H137").Value = 2037.1428
$ws.Range("I137").Value = 1520.0667
$ws.Range("J137").Value = 3329.8333
$ws.Range("K137").Value = 4560.2001
$ws.Range("L137").Value = 9989.499899999999
$ws.Range("M137").Value = -2010.2001
$ws.Range("N137").Value = -15089.4999
$ws.Range("H138").Value = 1843.7222
$ws.Range("I138").Value = 1649.0555
$ws.Range("J138").Value = 2038.3889
$ws.Range("K138").Value = 4947.166499999999
$ws.Range("L138").Value = 6115.1667
$ws.Range("M138").Value = 192.8335000000006
$ws.Range("N138").Value = -16395.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4558.6724
$ws.Range("I32").Value = 3126.4807
$ws.Range("K32").Value = 3126.4807
$ws.Range("M32").Value = -2839.4807
$ws.Range("H45").Value = 1822.2307
$ws.Range("I45").Value = 1499
$ws.Range("J45").Value = 1881
$ws.Range("K45").Value = 1499
$ws.Range("L45").Value = 1881
$ws.Range("M45").Value = -1122
$ws.Range("N45").Value = -2635
$ws.Range("H110").Value = 132.57143
$ws.Range("I110").Value = 132.57143
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 132.57143
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1912.42857
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 443.45456
$ws.Range("I107").Value = 396.88235
$ws.Range("K107").Value = 396.88235
$ws.Range("M107").Value = 1523.11765
$ws.Range("H134").Value = 7007.3687
$ws.Range("I134").Value = 8730.929
$ws.Range("K134").Value = 26192.787
$ws.Range("M134").Value = -23657.787

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2999.4443
$ws.Range("I31").Value = 2999.8572
$ws.Range("K31").Value = 2999.8572
$ws.Range("M31").Value = -2704.8572
$ws.Range("H34").Value = 2999.4443
$ws.Range("I34").Value = 2999.8572
$ws.Range("K34").Value = 2999.8572
$ws.Range("M34").Value = -2797.8572
$ws.Range("H99").Value = 3222.2222
$ws.Range("J99").Value = 4017.6
$ws.Range("L99").Value = 4017.6
$ws.Range("N99").Value = -7013.6
$ws.Range("H126").Value = 3222.2222
$ws.Range("J126").Value = 4017.6
$ws.Range("L126").Value = 12052.8
$ws.Range("N126").Value = -16992.8
$ws.Range("H134").Value = 2057.2
$ws.Range("I134").Value = 1710.4546
$ws.Range("K134").Value = 5131.3638
$ws.Range("M134").Value = -2596.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11583.682
$ws.Range("J131").Value = 11746.046
$ws.Range("L131").Value = 35238.138
$ws.Range("N131").Value = -45318.138

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2792.6155
$ws.Range("I102").Value = 2664.2727
$ws.Range("J102").Value = 3498.5
$ws.Range("K102").Value = 2664.2727
$ws.Range("L102").Value = 3498.5
$ws.Range("M102").Value = -1042.2727
$ws.Range("N102").Value = -6742.5
$ws.Range("H113").Value = 1104.091
$ws.Range("I113").Value = 960.75
$ws.Range("J113").Value = 1186
$ws.Range("K113").Value = 960.75
$ws.Range("L113").Value = 1186
$ws.Range("M113").Value = 1209.25
$ws.Range("N113").Value = -5526
$ws.Range("H122").Value = 1570.6923
$ws.Range("I122").Value = 1434.4445
$ws.Range("K122").Value = 4303.333500000001
$ws.Range("M122").Value = -1853.333500000001
$ws.Range("H132").Value = 1834834.5
$ws.Range("I132").Value = 2567082.2
$ws.Range("K132").Value = 7701246.600000001
$ws.Range("M132").Value = -7698716.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2109.4546
$ws.Range("I7").Value = 1817.1875
$ws.Range("K7").Value = 1817.1875
$ws.Range("M7").Value = -1705.1875
$ws.Range("H126").Value = 2109.4546
$ws.Range("I126").Value = 1817.1875
$ws.Range("K126").Value = 5451.5625
$ws.Range("M126").Value = -2981.5625
$ws.Range("H131").Value = 46130.4
$ws.Range("J131").Value = 46130.4
$ws.Range("L131").Value = 46130.4
$ws.Range("N131").Value = -56210.4
$ws.Range("H132").Value = 2303.2273
$ws.Range("I132").Value = 1284.5
$ws.Range("K132").Value = 3853.5
$ws.Range("M132").Value = -1323.5

